$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 9957446
$ws.Range("J32").Value = 11616853
$ws.Range("L32").Value = 11616853
$ws.Range("N32").Value = -11617505

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 3094.8333
$ws.Range("I86").Value = 2283.5
$ws.Range("J86").Value = 4022.0715
$ws.Range("K86").Value = 2283.5
$ws.Range("L86").Value = 4022.0715
$ws.Range("M86").Value = -1160.5
$ws.Range("N86").Value = -6268.0715

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 3094.8333
$ws.Range("I89").Value = 2283.5
$ws.Range("J89").Value = 4022.0715
$ws.Range("K89").Value = 11417.5
$ws.Range("L89").Value = 20110.3575
$ws.Range("M89").Value = -5801.5
$ws.Range("N89").Value = -31342.3575

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 32474.166
$ws.Range("I116").Value = 46974.25
$ws.Range("K116").Value = 46974.25
$ws.Range("M116").Value = -43532.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2258.5781
$ws.Range("I132").Value = 1254.2931
$ws.Range("J132").Value = 11966.667
$ws.Range("K132").Value = 3762.879300000001
$ws.Range("L132").Value = 35900.001
$ws.Range("M132").Value = -1232.879300000001
$ws.Range("N132").Value = -40960.001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 4267.2173
$ws.Range("I122").Value = 4375.6924
$ws.Range("J122").Value = 4126.2
$ws.Range("K122").Value = 13127.0772
$ws.Range("L122").Value = 12378.6
$ws.Range("M122").Value = -10677.0772
$ws.Range("N122").Value = -17278.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3279.9565
$ws.Range("I99").Value = 4393.2666
$ws.Range("J99").Value = 1192.5
$ws.Range("K99").Value = 4393.2666
$ws.Range("L99").Value = 1192.5
$ws.Range("M99").Value = -2895.2666
$ws.Range("N99").Value = -4188.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 971.36365
$ws.Range("I16").Value = 965.55554
$ws.Range("J16").Value = 997.5
$ws.Range("K16").Value = 965.55554
$ws.Range("L16").Value = 997.5
$ws.Range("M16").Value = -678.55554
$ws.Range("N16").Value = -1571.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2247.0405
$ws.Range("I31").Value = 1720.4286
$ws.Range("J31").Value = 3279.2
$ws.Range("K31").Value = 1720.4286
$ws.Range("L31").Value = 3279.2
$ws.Range("M31").Value = -1425.4286
$ws.Range("N31").Value = -3869.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2247.0405
$ws.Range("I34").Value = 1720.4286
$ws.Range("J34").Value = 3279.2
$ws.Range("K34").Value = 1720.4286
$ws.Range("L34").Value = 3279.2
$ws.Range("M34").Value = -1518.4286
$ws.Range("N34").Value = -3683.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4611.436
$ws.Range("I58").Value = 6702.6113
$ws.Range("J58").Value = 2819
$ws.Range("K58").Value = 6702.6113
$ws.Range("L58").Value = 2819
$ws.Range("M58").Value = -6499.6113
$ws.Range("N58").Value = -3225

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1197.1428
$ws.Range("I105").Value = 1037.5
$ws.Range("J105").Value = 1410
$ws.Range("K105").Value = 1037.5
$ws.Range("L105").Value = 1410
$ws.Range("M105").Value = 709.5
$ws.Range("N105").Value = -4904

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 228.625
$ws.Range("I107").Value = 171.89473
$ws.Range("J107").Value = 444.2
$ws.Range("K107").Value = 171.89473
$ws.Range("L107").Value = 444.2
$ws.Range("M107").Value = 1748.10527
$ws.Range("N107").Value = -4284.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 971.36365
$ws.Range("I113").Value = 965.55554
$ws.Range("J113").Value = 997.5
$ws.Range("K113").Value = 965.55554
$ws.Range("L113").Value = 997.5
$ws.Range("M113").Value = 1204.44446
$ws.Range("N113").Value = -5337.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1391.1063
$ws.Range("I134").Value = 921.7714
$ws.Range("J134").Value = 2760
$ws.Range("K134").Value = 2765.3142
$ws.Range("L134").Value = 8280
$ws.Range("M134").Value = -230.3141999999998
$ws.Range("N134").Value = -13350

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H135").Value = 38635.855
$ws.Range("J135").Value = 26912.666
$ws.Range("L135").Value = 26912.666
$ws.Range("N135").Value = -37052.666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 4611.436
$ws.Range("I136").Value = 6702.6113
$ws.Range("J136").Value = 2819
$ws.Range("K136").Value = 20107.8339
$ws.Range("L136").Value = 8457
$ws.Range("M136").Value = -17557.8339
$ws.Range("N136").Value = -13557

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3135.5
$ws.Range("I126").Value = 2880
$ws.Range("J126").Value = 3305.8333
$ws.Range("K126").Value = 8640
$ws.Range("L126").Value = 9917.499899999999
$ws.Range("M126").Value = -6170
$ws.Range("N126").Value = -14857.4999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3506.4036
$ws.Range("I136").Value = 1885.3572
$ws.Range("J136").Value = 8045.3335
$ws.Range("K136").Value = 5656.071599999999
$ws.Range("L136").Value = 24136.0005
$ws.Range("M136").Value = -3106.071599999999
$ws.Range("N136").Value = -29236.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2316.6667
$ws.Range("I81").Value = 1192.8572
$ws.Range("J81").Value = 6250
$ws.Range("K81").Value = 2385.7144
$ws.Range("L81").Value = 12500
$ws.Range("M81").Value = -1324.7144
$ws.Range("N81").Value = -14622

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 2316.6667
$ws.Range("I84").Value = 1192.8572
$ws.Range("J84").Value = 6250
$ws.Range("K84").Value = 11928.572
$ws.Range("L84").Value = 62500
$ws.Range("M84").Value = -6624.572
$ws.Range("N84").Value = -73108

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H88").Value = 39000
$ws.Range("J88").Value = 39000
$ws.Range("L88").Value = 39000
$ws.Range("N88").Value = -39812

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H91").Value = 39000
$ws.Range("J91").Value = 39000
$ws.Range("L91").Value = 39000
$ws.Range("N91").Value = -41808

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 27995
$ws.Range("J92").Value = 27995
$ws.Range("L92").Value = 27995
$ws.Range("N92").Value = -32987

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H106").Value = 25675
$ws.Range("J106").Value = 25675
$ws.Range("L106").Value = 25675
$ws.Range("N106").Value = -28199

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H124").Value = 20999.5
$ws.Range("J124").Value = 20999.5
$ws.Range("L124").Value = 20999.5
$ws.Range("N124").Value = -30819.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 875.43243
$ws.Range("I126").Value = 464.47058
$ws.Range("J126").Value = 1224.75
$ws.Range("K126").Value = 1393.41174
$ws.Range("L126").Value = 3674.25
$ws.Range("M126").Value = 1076.58826
$ws.Range("N126").Value = -8614.25
